# Updates the cryptos list data (price / 1h volume columns, plus a
# reordering of the Stacks/ImmutableX rows) to match the refreshed
# GitHub Actions data pull.
#
# Cells whose new text looks like a plain number (e.g. "605.57",
# "2.00", "0.999") are force-formatted as Text ("@") before the value
# is written so they are stored as literal strings (matching the
# original inline-string cells) instead of being auto-coerced to
# numeric values by Excel's normal type inference.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.169.14"
$ws.Range("E2").Value = "  +3.10%  "
$ws.Range("D3").Value = "2.618.97"
$ws.Range("E3").Value = "  +3.60%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.57"
$ws.Range("E5").Value = "  +1.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.21"
$ws.Range("E6").Value = "  +1.23%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +0.79%  "
$ws.Range("D9").Value = "2.617.14"
$ws.Range("E9").Value = "  +3.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.167"
$ws.Range("E10").Value = "  +12.79%  "
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.348"
$ws.Range("E12").Value = "  +2.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.03"
$ws.Range("E13").Value = "  +0.81%  "
$ws.Range("D14").Value = "3.097.06"
$ws.Range("E14").Value = "  +4.60%  "
$ws.Range("E15").Value = "  +7.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.68"
$ws.Range("E16").Value = "  +1.54%  "
$ws.Range("D17").Value = "71.116.02"
$ws.Range("E17").Value = "  +3.05%  "
$ws.Range("D18").Value = "2.614.53"
$ws.Range("E18").Value = "  +5.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "381.58"
$ws.Range("E19").Value = "  +6.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.48"
$ws.Range("E20").Value = "  +2.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.87"
$ws.Range("E21").Value = "  +3.82%  "
$ws.Range("E22").Value = "  -0.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.00"
$ws.Range("E23").Value = "  +17.91%  "
$ws.Range("E24").Value = "  +2.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.45"
$ws.Range("E25").Value = "  +5.66%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.08"
$ws.Range("E27").Value = "  +12.26%  "
$ws.Range("D28").Value = "2.752.83"
$ws.Range("E28").Value = "  +3.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "548.88"
$ws.Range("E30").Value = "  +4.94%  "
$ws.Range("D31").Value = "0.0₃0955"
$ws.Range("E31").Value = "  +7.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.05"
$ws.Range("E32").Value = "  +3.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.34"
$ws.Range("E33").Value = "  +6.84%  "
$ws.Range("E34").Value = "  +2.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "164.88"
$ws.Range("E36").Value = "  +0.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.115"
$ws.Range("E37").Value = "  -4.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.20"
$ws.Range("E38").Value = "  +4.00%  "
$ws.Range("E39").Value = "  +2.45%  "
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.41"
$ws.Range("E40").Value = "  +6.30%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.88"
$ws.Range("E41").Value = "  +5.44%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("E43").Value = "  +9.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.02"
$ws.Range("E44").Value = "  +3.53%  "
$ws.Range("E45").Value = "  +1.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.96"
$ws.Range("E46").Value = "  +2.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "153.04"
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.63"
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("E49").Value = "  +2.57%  "
$ws.Range("E50").Value = "  +4.36%  "
$ws.Range("E51").Value = "  +3.36%  "
